$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("D6").Formula = "=2*0.353+2*0.359+0.46"

# Row 7
$ws.Range("D7").Formula = "=0.646*2"

# Row 8
$ws.Range("D8").Formula = "=0.404*2"

# Row 9
$ws.Range("D9").Formula = "=2*(0.37+0.379+0.194+0.32+0.274+0.116+0.112+0.272)"

# Row 11 (hard-coded value)
$ws.Range("D11").Value = 2.047

# Row 12
$ws.Range("D12").Formula = "=0.22+0.218+0.116+0.057+0.056"

# Row 13
$ws.Range("D13").Formula = "=2*(0.768+0.708+0.258+0.546+0.404)"

# Row 14
$ws.Range("D14").Formula = "=0.537+0.233+0.292"

# Row 15 (hard-coded value)
$ws.Range("D15").Value = 0.162

# Row 16 (hard-coded value)
$ws.Range("D16").Value = 3.573

# Row 17
$ws.Range("D17").Formula = "=2*0.779"

# Row 18
$ws.Range("D18").Formula = "=2*(0.22+0.222+0.289+0.37+0.199+0.478)+0.446"

# Row 19 (hard-coded value)
$ws.Range("D19").Value = 0.707

# Update the active selection to match the latest edit location
$ws.Range("D21").Select()
